$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column B: 'Keyed' class label for all new rows, matching style of rows 160-167 ---
$ws.Range("B168:B175").Value = "Keyed"
$ws.Range("B168:B175").Interior.Color = 16777215

# --- Column F: English source text w/ file reference, written first (row order) ---
$ws.Range("F168").Value = 'MTP.OnlyWhenDraftedTT ''Will only apply when pawns are drafted'' (English file: Text.xml:62)'
$ws.Range("F169").Value = 'MTP.OrLogic ''Or-logic'' (English file: Text.xml:63)'
$ws.Range("F170").Value = 'MTP.OrLogicTT ''Will apply the rule if any of the selected things are correct, instead of when all are correct.'' (English file: Text.xml:64)'
$ws.Range("F171").Value = 'MTP.AutomaticType.TDFindLib ''Advanced rule (TD Find Lib)'' (English file: Text.xml:141)'
$ws.Range("F172").Value = 'MTP.AutomaticType.Mechanoid ''Specific mechanoid'' (English file: Text.xml:142)'
$ws.Range("F173").Value = 'MTP.TDFindLibRuleLabel ''Mark That Pawn-rule'' (English file: Text.xml:143)'
$ws.Range("F174").Value = 'MTP.RequiresAnActiveGame ''Requires an active game'' (English file: Text.xml:144)'
$ws.Range("F175").Value = 'MTP.EditTdRule ''Edit'' (English file: Text.xml:145)'

# --- Column G/H: parsing formulas (shared across 169:175, standalone on 168) ---
$ws.Range("G168").Formula = '=LEFT(F168,FIND(" ",F168)-1)'
$ws.Range("H168").Formula = '=MID(F168,FIND("''",F168)+1,FIND("''",F168,FIND("''",F168)+1)-FIND("''",F168)-1)'
$ws.Range("G169:G175").Formula = '=LEFT(F169,FIND(" ",F169)-1)'
$ws.Range("H169:H175").Formula = '=MID(F169,FIND("''",F169)+1,FIND("''",F169,FIND("''",F169)+1)-FIND("''",F169)-1)'

# --- Column C: key names. Written in this exact row order to reproduce the shared-
#     string table layout of the original authored workbook (row 169's key was
#     typed last in this column during the original edit). ---
$ws.Range("C168").Value = 'MTP.OnlyWhenDraftedTT'
$ws.Range("C170").Value = 'MTP.OrLogicTT'
$ws.Range("C171").Value = 'MTP.AutomaticType.TDFindLib'
$ws.Range("C172").Value = 'MTP.AutomaticType.Mechanoid'
$ws.Range("C173").Value = 'MTP.TDFindLibRuleLabel'
$ws.Range("C174").Value = 'MTP.RequiresAnActiveGame'
$ws.Range("C175").Value = 'MTP.EditTdRule'
$ws.Range("C169").Value = 'MTP.OrLogic'

# --- Column D: English values, rows 168-174 (row 175 'Edit' is typed last, after E) ---
$ws.Range("D168").Value = 'Will only apply when pawns are drafted'
$ws.Range("D169").Value = 'Or-logic'
$ws.Range("D170").Value = 'Will apply the rule if any of the selected things are correct, instead of when all are correct.'
$ws.Range("D171").Value = 'Advanced rule (TD Find Lib)'
$ws.Range("D172").Value = 'Specific mechanoid'
$ws.Range("D173").Value = 'Mark That Pawn-rule'
$ws.Range("D174").Value = 'Requires an active game'

# --- Column E: Korean translations, rows 168-174 ---
$ws.Range("E168").Value = '폰을 소집할 때만 적용됩니다.'
$ws.Range("E169").Value = 'Or-논리'
$ws.Range("E170").Value = '모든 항목이 맞을 때가 아니라 선택한 항목 중 하나라도 맞으면 규칙을 적용합니다.'
$ws.Range("E171").Value = '고급 규칙(TD Find Lib)'
$ws.Range("E172").Value = '특정 메카노이드'
$ws.Range("E173").Value = 'Mark That Pawn-규칙'
$ws.Range("E174").Value = '게임이 활성화되어 있어야 합니다.'

# --- Row 175 tail: KO translation reuses the existing '편집' string, then EN 'Edit' ---
$ws.Range("E175").Value = '편집'
$ws.Range("D175").Value = 'Edit'

# --- Update the active selection to match the authored workbook ---
$ws.Range("A162").Select()

